$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) stores these updated values as text,
# matching the existing inline-string cell type used throughout the sheet.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.152.54"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "1.672.11"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.79%  "
$ws.Range("D5").Value = "210.72"
$ws.Range("E5").Value = "  -3.87%  "
$ws.Range("D6").Value = "0.5253"
$ws.Range("E6").Value = "  -4.50%  "
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").Value = "0.2653"
$ws.Range("E8").Value = "  -3.17%  "
$ws.Range("D9").Value = "0.06288"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("D10").Value = "21.21"
$ws.Range("E10").Value = "  -3.68%  "
$ws.Range("D11").Value = "0.07518"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "1.692.58"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "4.444"
$ws.Range("E13").Value = "  -2.17%  "
$ws.Range("D14").Value = "0.5633"
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").Value = "0.000008029"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("D16").Value = "66.52"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "26.205.88"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "4.802"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").Value = "187.80"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").Value = "10.38"
$ws.Range("E21").Value = "  -5.46%  "
$ws.Range("D22").Value = "6.180"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "0.1249"
$ws.Range("E25").Value = "  -5.92%  "
$ws.Range("D26").Value = "7.598"
$ws.Range("E26").Value = "  -4.03%  "
$ws.Range("D27").Value = "15.96"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("D28").Value = "0.06234"
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").Value = "1.278"
$ws.Range("E30").Value = "  -4.17%  "
$ws.Range("D31").Value = "3.476"
$ws.Range("E31").Value = "  -3.63%  "
$ws.Range("D32").Value = "3.435"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("D33").Value = "1.624"
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("D34").Value = "0.9960"
$ws.Range("E34").Value = "  -4.45%  "
$ws.Range("D35").Value = "0.6040"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").Value = "6.112"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").Value = "0.01613"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").Value = "1.074.78"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").Value = "0.8672"
$ws.Range("E41").Value = "  -2.09%  "
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("D43").Value = "100.03"
$ws.Range("E43").Value = "  -1.76%  "
$ws.Range("D44").Value = "1.821.40"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "56.08"
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("D47").Value = "0.9999"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("D48").Value = "0.05241"
$ws.Range("E48").Value = "  -0.88%  "
$ws.Range("D49").Value = "7.978"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").Value = "0.4254"
$ws.Range("D51").Value = "6.000"
$ws.Range("E51").Value = "  -1.87%  "
